$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the whole used range first, since we're shrinking the data from 17 rows to 2.
$ws.Cells.Clear()

# Header row
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Teléfono"
$ws.Range("C1").Value = "Servicio"
$ws.Range("D1").Value = "Día y Hora"

# Data row
$ws.Range("A2").Value = "silvanito"
$ws.Range("B2").Value = 927859435
$ws.Range("C2").Value = "internet_cable"
$ws.Range("D2").Value = "2024-07-07 20:30:23"
